$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '34.575.73'
$ws.Range('E2').Value = '  +1.27%  '

# Row 3
$ws.Range('D3').Value = '1.797.28'
$ws.Range('E3').Value = '  +1.02%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.28'
$ws.Range('E5').Value = '  +0.53%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.556'
$ws.Range('E6').Value = '  +1.92%  '

# Row 7
$ws.Range('E7').Value = '  -0.09%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '32.82'
$ws.Range('E8').Value = '  +3.50%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.295'
$ws.Range('E9').Value = '  +1.35%  '

# Row 10
$ws.Range('E10').Value = '  +0.39%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0950'
$ws.Range('E11').Value = '  +0.50%  '

# Row 12
$ws.Range('D12').Value = '2.058.96'
$ws.Range('E12').Value = '  +1.13%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.819.87'
$ws.Range('E13').Value = '  +2.48%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.10'
$ws.Range('E14').Value = '  +1.42%  '

# Row 15
$ws.Range('E15').Value = '  +2.77%  '

# Row 16
$ws.Range('D16').Value = '34.546.83'
$ws.Range('E16').Value = '  +1.28%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '4.31'
$ws.Range('E17').Value = '  +3.08%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '68.85'
$ws.Range('E18').Value = '  +1.48%  '

# Row 19
$ws.Range('D19').Value = '0.0₃0802'
$ws.Range('E19').Value = '  +0.70%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '246.56'
$ws.Range('E20').Value = '  +0.51%  '

# Row 21
$ws.Range('E21').Value = '  +3.56%  '

# Row 22
$ws.Range('E22').Value = '  -0.20%  '

# Row 23
$ws.Range('E23').Value = '  +1.74%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '174.66'
$ws.Range('E24').Value = '  +7.51%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.06'
$ws.Range('E25').Value = '  +0.89%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.32'
$ws.Range('E26').Value = '  +1.73%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.69'
$ws.Range('E27').Value = '  +2.62%  '

# Row 28
$ws.Range('E28').Value = '  +1.44%  '

# Row 29
$ws.Range('E29').Value = '  -0.02%  '

# Row 30
$ws.Range('E30').Value = '  +8.41%  '

# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.24'
$ws.Range('E31').Value = '  +0.58%  '

# Row 32
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0525'
$ws.Range('E32').Value = '  +1.00%  '

# Row 33
$ws.Range('E33').Value = '  +1.38%  '

# Row 34
$ws.Range('E34').Value = '  +2.10%  '

# Row 35
$ws.Range('D35').Value = '1.429.80'
$ws.Range('E35').Value = '  -0.97%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.56'
$ws.Range('E36').Value = '  +6.50%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.675'
$ws.Range('E37').Value = '  +2.18%  '

# Row 38
$ws.Range('E38').Value = '  +2.41%  '

# Row 39
$ws.Range('E39').Value = '  +0.53%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '84.74'
$ws.Range('E40').Value = '  +5.66%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.947'
$ws.Range('E41').Value = '  +2.90%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.39'
$ws.Range('E42').Value = '  +1.29%  '

# Row 43
$ws.Range('E43').Value = '  +3.17%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.90'
$ws.Range('E44').Value = '  +3.38%  '

# Row 45
$ws.Range('E45').Value = '  +2.73%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '6.09'
$ws.Range('E46').Value = '  +0.40%  '

# Row 47
$ws.Range('E47').Value = '  +1.15%  '

# Row 48
$ws.Range('E48').Value = '  +1.09%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '105.25'
$ws.Range('E49').Value = '  +1.26%  '

# Row 50
$ws.Range('E50').Value = '  -0.06%  '

# Row 51
$ws.Range('D51').Value = '0.0₆0129'
$ws.Range('E51').Value = '  -4.89%  '
